$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.352.24'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  +1.11%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.606.35'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  +0.92%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '541.81'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  +4.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.71'
$ws.Range("D6").Style = "Normal"

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("E8").Value = '  +0.33%  '

$ws.Range("E9").Value = '  -1.59%  '

$ws.Range("E10").Value = '  +2.31%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.337'
$ws.Range("D11").Style = "Normal"

$ws.Range("E12").Value = '  +0.54%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.060.06'
$ws.Range("D13").Style = "Normal"

$ws.Range("E13").Value = '  +0.78%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '59.285.34'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  +1.03%  '

$ws.Range("E15").Value = '  +1.34%  '

$ws.Range("E16").Value = '  +1.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.605.48'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  -0.65%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '341.35'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  +0.91%  '

$ws.Range("E19").Value = '  +1.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.14'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  +0.17%  '

$ws.Range("E21").Value = '  -1.38%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"

$ws.Range("E22").Value = '  +0.07%  '

$ws.Range("E24").Value = '  +1.58%  '

$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.26'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  +3.42%  '

$ws.Range("E28").Value = '  +4.03%  '

$ws.Range("E30").Value = '  +6.97%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.82'
$ws.Range("D31").Style = "Normal"

$ws.Range("E31").Value = '  -1.73%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.77'
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '149.87'
$ws.Range("D33").Style = "Normal"

$ws.Range("E33").Value = '  +1.15%  '

$ws.Range("E34").Value = '  +0.42%  '

$ws.Range("E35").Value = '  +0.20%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '37.23'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  +2.09%  '

$ws.Range("E37").Value = '  +0.74%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.838'
$ws.Range("D38").Style = "Normal"

$ws.Range("E38").Value = '  +1.27%  '

$ws.Range("E39").Value = '  +1.83%  '

$ws.Range("E40").Value = '  +1.89%  '

$ws.Range("E41").Value = '  +0.17%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '276.01'
$ws.Range("D42").Style = "Normal"

$ws.Range("E42").Value = '  +0.60%  '

$ws.Range("E43").Value = '  +2.05%  '

$ws.Range("E44").Value = '  -0.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0956'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  +0.89%  '

$ws.Range("E46").Value = '  +1.20%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.953.00'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  -1.12%  '

$ws.Range("B48").Value = 'VeChain'

$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0224'
$ws.Range("D48").Style = "Normal"

$ws.Range("E48").Value = '  +2.09%  '

$ws.Range("B49").Value = 'InjectiveProtocol'

$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.56'
$ws.Range("D49").Style = "Normal"

$ws.Range("E49").Value = '  +4.00%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.53'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  +0.62%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '110.87'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  -1.07%  '
